$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D-column values must remain literal text (e.g. "1.000" must not become 1),
# so we force a temporary text format, assign the value, then restore the
# original cell style so formatting/style stays unchanged.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.411.30"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +0.73%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.639.32"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +2.29%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  -0.03%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -0.01%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "304.53"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +0.39%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3734"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -1.01%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "52.13"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +0.47%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3622"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -0.31%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.247"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -2.31%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08101"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -0.35%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -0.01%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.81"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -0.10%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.588"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -0.33%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001268"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +1.48%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.278"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -2.01%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.630.04"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +1.56%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "94.37"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +0.41%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06887"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -0.56%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.11"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.26%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.512"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -0.35%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +0.05%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.421.39"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +0.74%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.75"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -1.70%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.411"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +1.21%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.037"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +0.15%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.19"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.18%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "151.67"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +1.00%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.340"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +1.55%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "135.50"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +0.79%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.283"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -4.80%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.808.14"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +1.50%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.781"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +0.13%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9493"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.49%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02819"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +2.69%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.31"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -0.35%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2522"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -0.72%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07208"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -4.15%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08759"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -0.86%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.074"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -0.94%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.369"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -1.90%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7028"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -1.24%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.43"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -0.70%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.96"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  +1.59%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6498"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -0.70%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.328"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +0.33%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9993"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -0.03%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.004"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -0.29%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07968"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.22%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "128.21"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -3.52%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.197"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -0.73%  "

